# Auto-generated Excel COM-interop script applying the scraped diff
# to Sheets/Chocobo_Profits.xlsx (workbook with sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Each block updates the commodity-profit columns (H:N) for specific rows,
# matching values refreshed by the scheduled market-data runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 101
$ws.Range("H101").Value = 1805.7273
$ws.Range("I101").Value = 993.3333
$ws.Range("K101").Value = 2979.9999
$ws.Range("M101").Value = -1357.9999
# Row 117
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
# Row 138
$ws.Range("H138").Value = 3274.0154
$ws.Range("J138").Value = 3914.739
$ws.Range("L138").Value = 11744.217
$ws.Range("N138").Value = -22024.217
# Row 141
$ws.Range("H141").Value = 8929.0625
$ws.Range("I141").Value = 9491.786
$ws.Range("K141").Value = 28475.358
$ws.Range("M141").Value = -23295.358

$ws = $wb.Worksheets.Item("ARM")
# Row 28
$ws.Range("H28").Value = 5604.25
$ws.Range("I28").Value = 5604.25
$ws.Range("K28").Value = 5604.25
$ws.Range("M28").Value = -5412.25
# Row 45
$ws.Range("H45").Value = 2011.1818
$ws.Range("I45").Value = 2224.125
$ws.Range("J45").Value = 1443.3334
$ws.Range("K45").Value = 2224.125
$ws.Range("L45").Value = 1443.3334
$ws.Range("M45").Value = -1847.125
$ws.Range("N45").Value = -2197.3334
# Row 61
$ws.Range("H61").Value = 1545.5652
$ws.Range("I61").Value = 1449.4736
$ws.Range("J61").Value = 2002
$ws.Range("K61").Value = 1449.4736
$ws.Range("L61").Value = 2002
$ws.Range("M61").Value = -1237.4736
$ws.Range("N61").Value = -2426
# Row 99
$ws.Range("H99").Value = 5604.25
$ws.Range("I99").Value = 5604.25
$ws.Range("K99").Value = 5604.25
$ws.Range("M99").Value = -2609.25
# Row 136
$ws.Range("H136").Value = 1545.5652
$ws.Range("I136").Value = 1449.4736
$ws.Range("J136").Value = 2002
$ws.Range("K136").Value = 4348.4208
$ws.Range("L136").Value = 6006
$ws.Range("M136").Value = -1798.4208
$ws.Range("N136").Value = -11106

$ws = $wb.Worksheets.Item("BSM")
# Row 7
$ws.Range("H7").Value = 6018778.5
$ws.Range("I7").Value = 6680000.5
$ws.Range("J7").Value = 5735397.5
$ws.Range("K7").Value = 6680000.5
$ws.Range("L7").Value = 5735397.5
$ws.Range("M7").Value = -6679887.5
$ws.Range("N7").Value = -5735623.5
# Row 107
$ws.Range("H107").Value = 2000
$ws.Range("I107").Value = 2000
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 2000
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = -80
$ws.Range("N107").Value = -5840
# Row 134
$ws.Range("H134").Value = 3351.1555
$ws.Range("I134").Value = 1201.2142
$ws.Range("J134").Value = 6892.2354
$ws.Range("K134").Value = 3603.6426
$ws.Range("L134").Value = 20676.7062
$ws.Range("M134").Value = -1068.6426
$ws.Range("N134").Value = -25746.7062

$ws = $wb.Worksheets.Item("CUL")
# Row 23
$ws.Range("H23").Value = 187.875
$ws.Range("I23").Value = 75
$ws.Range("J23").Value = 225.5
$ws.Range("K23").Value = 225
$ws.Range("L23").Value = 676.5
$ws.Range("M23").Value = 10
$ws.Range("N23").Value = -1146.5
# Row 68
$ws.Range("H68").Value = 1342.2428
$ws.Range("I68").Value = 876.0294
$ws.Range("K68").Value = 2628.0882
$ws.Range("M68").Value = -1817.0882
# Row 71
$ws.Range("H71").Value = 1342.2428
$ws.Range("I71").Value = 876.0294
$ws.Range("K71").Value = 7884.2646
$ws.Range("M71").Value = -3828.2646
# Row 96
$ws.Range("H96").Value = 333336670
$ws.Range("I96").Value = 1000000000
$ws.Range("J96").Value = 5000
$ws.Range("K96").Value = 3000000000
$ws.Range("L96").Value = 15000
$ws.Range("M96").Value = -2999997941
$ws.Range("N96").Value = -19118
# Row 107
$ws.Range("H107").Value = 6186523.5
$ws.Range("J107").Value = 9454590
$ws.Range("L107").Value = 28363770
$ws.Range("N107").Value = -28367610

$ws = $wb.Worksheets.Item("GSM")
# Row 23
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
# Row 43
$ws.Range("H43").Value = 19583.37
$ws.Range("I43").Value = 1416.8334
$ws.Range("J43").Value = 27967.924
$ws.Range("K43").Value = 1416.8334
$ws.Range("L43").Value = 27967.924
$ws.Range("M43").Value = -1265.8334
$ws.Range("N43").Value = -28269.924
# Row 57
$ws.Range("H57").Value = 17817.908
$ws.Range("J57").Value = 17817.908
$ws.Range("L57").Value = 17817.908
$ws.Range("N57").Value = -19457.908
# Row 70
$ws.Range("H70").Value = 6820.2646
$ws.Range("I70").Value = 6049.231
$ws.Range("J70").Value = 9326.125
$ws.Range("K70").Value = 6049.231
$ws.Range("L70").Value = 9326.125
$ws.Range("M70").Value = -5779.231
$ws.Range("N70").Value = -9866.125
# Row 73
$ws.Range("H73").Value = 6820.2646
$ws.Range("I73").Value = 6049.231
$ws.Range("J73").Value = 9326.125
$ws.Range("K73").Value = 6049.231
$ws.Range("L73").Value = 9326.125
$ws.Range("M73").Value = -5113.231
# Row 80
$ws.Range("H80").Value = 15628187
$ws.Range("I80").Value = 22730508
$ws.Range("J80").Value = 3081.2
$ws.Range("K80").Value = 22730508
$ws.Range("L80").Value = 3081.2
$ws.Range("M80").Value = -22729510
$ws.Range("N80").Value = -5077.2
# Row 83
$ws.Range("H83").Value = 15628187
$ws.Range("I83").Value = 22730508
$ws.Range("J83").Value = 3081.2
$ws.Range("K83").Value = 113652540
$ws.Range("L83").Value = 15406
$ws.Range("M83").Value = -113647548
$ws.Range("N83").Value = -25390
# Row 94
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 16668247
$ws.Range("I22").Value = 27778888
$ws.Range("J22").Value = 2283.3333
$ws.Range("K22").Value = 27778888
$ws.Range("L22").Value = 2283.3333
$ws.Range("M22").Value = -27778593
$ws.Range("N22").Value = -2873.3333
# Row 27
$ws.Range("H27").Value = 16668247
$ws.Range("I27").Value = 27778888
$ws.Range("J27").Value = 2283.3333
$ws.Range("K27").Value = 27778888
$ws.Range("L27").Value = 2283.3333
$ws.Range("M27").Value = -27778781
# Row 69
$ws.Range("H69").Value = 40000
$ws.Range("J69").Value = 40000
$ws.Range("L69").Value = 40000
$ws.Range("N69").Value = -41622
# Row 72
$ws.Range("H72").Value = 40000
$ws.Range("J72").Value = 40000
$ws.Range("L72").Value = 120000
$ws.Range("N72").Value = -128112

$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 64020612
$ws.Range("I96").Value = 90955450
$ws.Range("J96").Value = 4763985
$ws.Range("K96").Value = 90955450
$ws.Range("L96").Value = 4763985
$ws.Range("M96").Value = -90954077
$ws.Range("N96").Value = -4766731
# Row 137
$ws.Range("H137").Value = 45131.11
$ws.Range("J137").Value = 45131.11
$ws.Range("L137").Value = 45131.11
$ws.Range("N137").Value = -55331.11
# Row 139
$ws.Range("H139").Value = 47508
$ws.Range("J139").Value = 47508
$ws.Range("L139").Value = 47508
$ws.Range("N139").Value = -57788
